# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# The underlying source data got re-ordered: several rows keep their row
# number / "id" sequence column (A) but the match data in columns B..AD
# (id/match number, teams, odds, etc.) moved to a different row.
#
# Net effect, expressed as "new row's B:AD content = old row's B:AD content":
#   row 14  <- old row 15
#   row 15  <- old row 14
#   row 427 <- old row 430
#   row 428 <- old row 431
#   row 429 <- old row 428
#   row 430 <- old row 427
#   row 431 <- old row 429
#
# We apply this by capturing the B:AD values of every row involved before
# any writes happen, then writing them into their new locations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($rowNumber) {
    return $ws.Range("B" + $rowNumber + ":AD" + $rowNumber).Value()
}

# Snapshot all rows that participate in the reshuffle before writing anything,
# since some rows both receive and donate data (3-cycle on 427/428/429/430/431).
$orig14  = Get-RowValues 14
$orig15  = Get-RowValues 15
$orig427 = Get-RowValues 427
$orig428 = Get-RowValues 428
$orig429 = Get-RowValues 429
$orig430 = Get-RowValues 430
$orig431 = Get-RowValues 431

# Pairwise swap: 14 <-> 15
$ws.Range("B14:AD14").Value = $orig15
$ws.Range("B15:AD15").Value = $orig14

# Pairwise swap: 427 <-> 430
$ws.Range("B427:AD427").Value = $orig430
$ws.Range("B430:AD430").Value = $orig427

# 3-cycle: 431 -> 428 -> 429 -> 431
$ws.Range("B428:AD428").Value = $orig431
$ws.Range("B429:AD429").Value = $orig428
$ws.Range("B431:AD431").Value = $orig429
